$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target "speaker_variant" (column C) values per row, after re-export
# with de-duplication / no preference flag / no levenshtein distance.
$values = @{
    2  = "eer"
    3  = "fav"
    4  = "timb"
    5  = "fen"
    6  = "bod"
    7  = "pag"
    8  = "tim"
    9  = "twe"
    10 = "lin"
    11 = "SIOYCK"
    12 = "SIOVCK"
    13 = "siouck"
    14 = "bal"
    15 = "Roemer"
    16 = "lio"
    17 = "1. VR"
    18 = "2. VR"
    19 = "alb"
    20 = "rod"
    21 = "siovck"
    22 = "gir"
    23 = "1, VR"
    24 = "doc"
    25 = "ROEMER"
    26 = "eph"
    27 = "ti"
    28 = "lyd"
    29 = "roemer"
    30 = "ped"
    31 = "leo"
}

# Rows 2-12 carried an "x" in is_prefered (column D) in the old export;
# the new export drops that flag entirely, so only those rows need D cleared.
$rowsWithPref = 2..12

foreach ($row in $values.Keys) {
    $variant = $values[$row]
    $id = "#" + ($variant.ToLower() -replace ", ", ",-" -replace "\. ", ".-" -replace " ", "-")

    $ws.Cells.Item($row, 2).Value = $id
    $ws.Cells.Item($row, 3).Value = $variant
}

foreach ($row in $rowsWithPref) {
    $ws.Cells.Item($row, 4).ClearContents()
}
